# Auto-generated edit script applying the diff's numeric cell updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 431.51163
$ws.Range("J17").Value = 431.51163
$ws.Range("L17").Value = 1294.53489
$ws.Range("N17").Value = -1630.53489
$ws.Range("H28").Value = 725.55
$ws.Range("I28").Value = 606.25
$ws.Range("J28").Value = 1202.75
$ws.Range("K28").Value = 606.25
$ws.Range("L28").Value = 1202.75
$ws.Range("M28").Value = -121.25
$ws.Range("N28").Value = -2172.75
$ws.Range("H38").Value = 400.8
$ws.Range("I38").Value = 400.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1202.4
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -830.4000000000001
$ws.Range("H45").Value = 1625
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1308
$ws.Range("H58").Value = 949.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 1591.6666
$ws.Range("I62").Value = 1410.2
$ws.Range("K62").Value = 1410.2
$ws.Range("M62").Value = -786.2
$ws.Range("H65").Value = 1591.6666
$ws.Range("I65").Value = 1410.2
$ws.Range("K65").Value = 7051
$ws.Range("M65").Value = -3931
$ws.Range("H68").Value = 80255
$ws.Range("J68").Value = 80255
$ws.Range("L68").Value = 80255
$ws.Range("N68").Value = -81753
$ws.Range("H71").Value = 80255
$ws.Range("J71").Value = 80255
$ws.Range("L71").Value = 240765
$ws.Range("N71").Value = -248253
$ws.Range("H111").Value = 2527.4167
$ws.Range("I111").Value = 2393.6365
$ws.Range("K111").Value = 7180.9095
$ws.Range("M111").Value = -4113.9095
$ws.Range("H114").Value = 126000
$ws.Range("J114").Value = 126000
$ws.Range("L114").Value = 126000
$ws.Range("N114").Value = -134678
$ws.Range("H132").Value = 3289.1924
$ws.Range("I132").Value = 3300.76
$ws.Range("K132").Value = 9902.280000000001
$ws.Range("M132").Value = -7372.280000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2772.077
$ws.Range("I45").Value = 2669.75
$ws.Range("K45").Value = 2669.75
$ws.Range("M45").Value = -2292.75
$ws.Range("H74").Value = 6104412
$ws.Range("I74").Value = 8335048
$ws.Range("K74").Value = 8335048
$ws.Range("M74").Value = -8334174
$ws.Range("H77").Value = 6104412
$ws.Range("I77").Value = 8335048
$ws.Range("K77").Value = 41675240
$ws.Range("M77").Value = -41670872
$ws.Range("H102").Value = 11174.462
$ws.Range("I102").Value = 11726.9
$ws.Range("K102").Value = 11726.9
$ws.Range("M102").Value = -10104.9
$ws.Range("H110").Value = 2801.8333
$ws.Range("I110").Value = 2671.1875
$ws.Range("K110").Value = 2671.1875
$ws.Range("M110").Value = -626.1875
$ws.Range("H132").Value = 4894.364
$ws.Range("I132").Value = 3344.9167
$ws.Range("J132").Value = 10363
$ws.Range("K132").Value = 10034.7501
$ws.Range("L132").Value = 31089
$ws.Range("M132").Value = -7504.750100000001
$ws.Range("N132").Value = -36149

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5127.25
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 254.5
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 254.5
$ws.Range("M5").Value = -9887
$ws.Range("N5").Value = -480.5
$ws.Range("H86").Value = 2791.125
$ws.Range("I86").Value = 2674.5715
$ws.Range("K86").Value = 2674.5715
$ws.Range("M86").Value = -1551.5715
$ws.Range("H89").Value = 2791.125
$ws.Range("I89").Value = 2674.5715
$ws.Range("K89").Value = 13372.8575
$ws.Range("M89").Value = -7756.8575
$ws.Range("H114").Value = 43100
$ws.Range("J114").Value = 43100
$ws.Range("L114").Value = 43100
$ws.Range("N114").Value = -51778
$ws.Range("H120").Value = 55320
$ws.Range("J120").Value = 55320
$ws.Range("L120").Value = 55320
$ws.Range("N120").Value = -64996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 442637.7
$ws.Range("I31").Value = 2957.7273
$ws.Range("K31").Value = 2957.7273
$ws.Range("M31").Value = -2662.7273
$ws.Range("H34").Value = 442637.7
$ws.Range("I34").Value = 2957.7273
$ws.Range("K34").Value = 2957.7273
$ws.Range("M34").Value = -2755.7273
$ws.Range("H98").Value = 55993.332
$ws.Range("J98").Value = 55993.332
$ws.Range("L98").Value = 55993.332
$ws.Range("N98").Value = -60485.332
$ws.Range("H117").Value = 109000
$ws.Range("J117").Value = 109000
$ws.Range("L117").Value = 109000
$ws.Range("N117").Value = -118178
$ws.Range("H119").Value = 80120
$ws.Range("J119").Value = 80120
$ws.Range("L119").Value = 80120
$ws.Range("N119").Value = -89796
$ws.Range("H134").Value = 306557.03
$ws.Range("I134").Value = 401772
$ws.Range("K134").Value = 1205316
$ws.Range("M134").Value = -1202781
$ws.Range("H138").Value = 69879.8
$ws.Range("J138").Value = 83422.5
$ws.Range("L138").Value = 83422.5
$ws.Range("N138").Value = -93702.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6685.625
$ws.Range("I3").Value = 2604.111
$ws.Range("J3").Value = 11933.286
$ws.Range("K3").Value = 7812.333
$ws.Range("L3").Value = 35799.858
$ws.Range("M3").Value = -7700.333
$ws.Range("N3").Value = -36023.858
$ws.Range("H12").Value = 594078.9
$ws.Range("I12").Value = 133.85715
$ws.Range("J12").Value = 1056036.1
$ws.Range("K12").Value = 401.57145
$ws.Range("L12").Value = 3168108.3
$ws.Range("M12").Value = -228.57145
$ws.Range("N12").Value = -3168454.3
$ws.Range("H51").Value = 13572.546
$ws.Range("I51").Value = 7923.5
$ws.Range("J51").Value = 28636.666
$ws.Range("K51").Value = 23770.5
$ws.Range("L51").Value = 85909.99800000001
$ws.Range("M51").Value = -23310.5
$ws.Range("N51").Value = -86829.99800000001
$ws.Range("H114").Value = 2305.4167
$ws.Range("I114").Value = 2593
$ws.Range("K114").Value = 7779
$ws.Range("M114").Value = -4525
$ws.Range("H120").Value = 17500
$ws.Range("H125").Value = 15005.5
$ws.Range("J125").Value = 15005.5
$ws.Range("L125").Value = 45016.5
$ws.Range("N125").Value = -54856.5
$ws.Range("H131").Value = 6428.263
$ws.Range("J131").Value = 4200.4375
$ws.Range("L131").Value = 12601.3125
$ws.Range("N131").Value = -22681.3125
$ws.Range("H134").Value = 3767.5667
$ws.Range("I134").Value = 2322.6365
$ws.Range("J134").Value = 7741.125
$ws.Range("K134").Value = 6967.9095
$ws.Range("L134").Value = 23223.375
$ws.Range("M134").Value = -1897.9095
$ws.Range("N134").Value = -33363.375
$ws.Range("H139").Value = 2694.389
$ws.Range("J139").Value = 2666.6667
$ws.Range("L139").Value = 8000.000100000001
$ws.Range("N139").Value = -18280.0001
$ws.Range("H140").Value = 378902.62
$ws.Range("I140").Value = 503848.66
$ws.Range("K140").Value = 1511545.98
$ws.Range("M140").Value = -1506365.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 99951
$ws.Range("J62").Value = 99951
$ws.Range("L62").Value = 99951
$ws.Range("N62").Value = -101323
$ws.Range("H65").Value = 99951
$ws.Range("J65").Value = 99951
$ws.Range("L65").Value = 299853
$ws.Range("N65").Value = -306717
$ws.Range("H113").Value = 3744.4211
$ws.Range("I113").Value = 3579.6667
$ws.Range("K113").Value = 3579.6667
$ws.Range("M113").Value = -1409.6667
$ws.Range("H122").Value = 1849.3334
$ws.Range("I122").Value = 1849.3334
$ws.Range("K122").Value = 5548.0002
$ws.Range("M122").Value = -3098.0002
$ws.Range("H132").Value = 22729688
$ws.Range("I132").Value = 24392662
$ws.Range("K132").Value = 73177986
$ws.Range("M132").Value = -73175456

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 48000
$ws.Range("J63").Value = 48000
$ws.Range("L63").Value = 48000
$ws.Range("N63").Value = -49498
$ws.Range("H66").Value = 48000
$ws.Range("J66").Value = 48000
$ws.Range("L66").Value = 144000
$ws.Range("N66").Value = -151488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 96742.25
$ws.Range("J124").Value = 96742.25
$ws.Range("L124").Value = 96742.25
$ws.Range("N124").Value = -106562.25
$ws.Range("H125").Value = 89784
$ws.Range("J125").Value = 89784
$ws.Range("L125").Value = 89784
$ws.Range("N125").Value = -99624
$ws.Range("H126").Value = 2060.9092
$ws.Range("I126").Value = 1492
$ws.Range("K126").Value = 4476
$ws.Range("M126").Value = -2006
$ws.Range("H128").Value = 67860
$ws.Range("J128").Value = 67860
$ws.Range("L128").Value = 67860
$ws.Range("N128").Value = -77820
$ws.Range("H129").Value = 94792.5
$ws.Range("J129").Value = 94792.5
$ws.Range("L129").Value = 94792.5
$ws.Range("N129").Value = -104792.5
$ws.Range("H131").Value = 78897
$ws.Range("J131").Value = 78897
$ws.Range("L131").Value = 78897
$ws.Range("N131").Value = -88977
$ws.Range("H132").Value = 16603.4
$ws.Range("I132").Value = 2613.2144
$ws.Range("K132").Value = 7839.6432
$ws.Range("M132").Value = -5309.6432
